# Apply the edits described by the commit diff:
# 1. Rename "Sheet1" -> "96-well map"
# 2. Make "96-well map" the active/selected tab, with selection on C24
# 3. Remove tab-selected state from "Platemap 11pt" (previously the active tab),
#    leaving its last selection on C22

$wb = $excel.ActiveWorkbook

# 1. Rename Sheet1 to "96-well map"
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "96-well map"

# 2. Select a cell on "96-well map" and make it the active sheet/tab
$sheet1.Activate()
$sheet1.Range("C24").Select()

# 3. Update selection on "Platemap 11pt" sheet (it previously had the active tab,
#    but the active tab now moves to "96-well map")
$platemap = $wb.Worksheets.Item("Platemap 11pt")
$platemap.Range("C22").Select()

# Re-activate "96-well map" so it remains the active/visible tab when saved
$sheet1.Activate()
